# Regenerate the attribute header columns (E1:K1) in the new order.
# Original order: IV, LO, LTT, OB, OM, SV, SO
# New order:      OM, SO, IV, OB, SV, LTT, LO
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E1:K1").ClearContents()

$ws.Range("E1").Value = "OM"
$ws.Range("F1").Value = "SO"
$ws.Range("G1").Value = "IV"
$ws.Range("H1").Value = "OB"
$ws.Range("I1").Value = "SV"
$ws.Range("J1").Value = "LTT"
$ws.Range("K1").Value = "LO"

# Column widths follow the (bestFit) width of the new header text in each
# column - mirror the widths Excel would have recalculated.
$ws.Columns("E").ColumnWidth = 4.142857142857143
$ws.Columns("F").ColumnWidth = 3.4285714285714284
$ws.Columns("G").ColumnWidth = 2.5714285714285716
$ws.Columns("H").ColumnWidth = 3.5714285714285716
$ws.Columns("I").ColumnWidth = 3.0
$ws.Columns("J").ColumnWidth = 3.7142857142857144
$ws.Columns("K").ColumnWidth = 3.2857142857142856
